$wb = $excel.ActiveWorkbook

$ws_FlowSegment = $wb.Worksheets.Item("FlowSegment")
$ws_FlowSegment_map = @{
    "M3" = 1935
    "O4" = 50
    "M5" = 1440
    "M6" = 1440
    "M7" = 1440
    "M8" = 1440
    "M11" = 2835
    "M12" = 2835
    "M13" = 2835
    "M14" = 2835
    "M15" = 2835
    "M16" = 2835
    "M17" = 1440
    "M18" = 1440
    "M21" = 1440
    "M22" = 1440
    "M23" = 1440
    "M24" = 1440
    "M25" = 1440
    "M26" = 1440
    "M27" = 2835
    "M29" = 2835
    "M30" = 2835
    "M31" = 2835
    "M32" = 1440
    "M33" = 2835
    "M34" = 1440
    "O36" = 50
    "O37" = 50
    "O38" = 50
    "O39" = 50
    "M42" = 1440
    "M43" = 1440
    "M45" = 2835
    "M46" = 2835
    "M49" = 2835
    "M50" = 2835
    "M51" = 2835
    "M54" = 2835
    "M56" = 2835
}
foreach ($key in $ws_FlowSegment_map.Keys) {
    $ws_FlowSegment.Range($key).Value = $ws_FlowSegment_map[$key]
}

$ws_FlowFitting = $wb.Worksheets.Item("FlowFitting")
$ws_FlowFitting_map = @{
    "M52" = 1440
}
foreach ($key in $ws_FlowFitting_map.Keys) {
    $ws_FlowFitting.Range($key).Value = $ws_FlowFitting_map[$key]
}

$ws_BuildingElementProxy = $wb.Worksheets.Item("BuildingElementProxy")
$ws_BuildingElementProxy_map = @{
    "M4" = 1440
    "M5" = 2835
    "M7" = 2835
    "M11" = 1935
    "M12" = 2835
    "M13" = 2835
    "M22" = 2835
    "M23" = 2835
    "M24" = 2835
    "M25" = 2835
    "M26" = 2835
    "M27" = 2835
    "M28" = 2835
    "M29" = 2835
    "M30" = 2835
    "M31" = 2835
    "M32" = 2835
    "M33" = 2835
    "M34" = 2835
    "M35" = 2835
    "M36" = 2835
    "M41" = 2835
    "M42" = 2835
    "M43" = 2835
    "M44" = 2835
    "M45" = 2835
    "M46" = 2835
    "M47" = 1935
    "M48" = 1935
    "M49" = 1440
    "M50" = 1440
    "M51" = 1440
    "M52" = 1440
    "M53" = 1440
    "M54" = 1440
    "M55" = 1440
    "M56" = 1440
    "M57" = 1440
    "M58" = 1440
    "M59" = 1440
    "M60" = 1440
    "M61" = 1440
    "M62" = 1440
    "M63" = 1440
    "M64" = 1440
    "M65" = 1440
    "M77" = 900
    "M78" = 2835
    "M79" = 1935
    "M80" = 1440
    "M98" = 900
    "M99" = 900
    "M100" = 900
    "M101" = 900
    "M102" = 900
    "M103" = 900
    "M104" = 900
    "M105" = 900
    "M112" = 1935
    "M113" = 1935
    "M114" = 1935
    "M115" = 1935
    "M116" = 1440
    "M117" = 1440
    "M118" = 1440
    "M119" = 1440
    "M120" = 1935
    "M121" = 2835
    "M122" = 2835
    "M123" = 2835
}
foreach ($key in $ws_BuildingElementProxy_map.Keys) {
    $ws_BuildingElementProxy.Range($key).Value = $ws_BuildingElementProxy_map[$key]
}

$ws_FlowTerminal = $wb.Worksheets.Item("FlowTerminal")
$ws_FlowTerminal_map = @{
    "M3" = 2835
    "M4" = 1440
}
foreach ($key in $ws_FlowTerminal_map.Keys) {
    $ws_FlowTerminal.Range($key).Value = $ws_FlowTerminal_map[$key]
}

$ws_Wall = $wb.Worksheets.Item("Wall")
$ws_Wall_map = @{
    "M2" = 2835
}
foreach ($key in $ws_Wall_map.Keys) {
    $ws_Wall.Range($key).Value = $ws_Wall_map[$key]
}

$ws_FurnishingElement = $wb.Worksheets.Item("FurnishingElement")
$ws_FurnishingElement_map = @{
    "M2" = 2835
}
foreach ($key in $ws_FurnishingElement_map.Keys) {
    $ws_FurnishingElement.Range($key).Value = $ws_FurnishingElement_map[$key]
}
